# Conserto do erro com o rotulo da coluna 2050 nas tabelas e retirada das
# linhas com total das tabelas.

$wb = $excel.ActiveWorkbook

# --- Sheets 1-3: "2040" header -> fix the trailing "2050" header in E1,
#     then drop the trailing "Total" row (row 13). ---
$sheetsSimple = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)"
)
foreach ($name in $sheetsSimple) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("E1").Value = "2050"
    $ws.Rows("13").Delete()
}

# --- Sheet 4: headers are periods ("2031-2040"), so E1 becomes "2041-2050".
#     Also drop the trailing "Total" row (row 13). ---
$ws = $wb.Worksheets.Item("Potencia Incremental - SIN(MW)")
$ws.Range("E1").Value = "2041-2050"
$ws.Rows("13").Delete()

# --- Sheet 5: same header fix, but this sheet never had a "Total" row. ---
$ws = $wb.Worksheets.Item("Emissoes Totais (MtCO2eq)")
$ws.Range("E1").Value = "2050"

# --- Sheet 6: no header in E1 to fix (only one data column); just drop the
#     trailing "Total" row (row 4). ---
$ws = $wb.Worksheets.Item("Custo Total (bilhões de R$)")
$ws.Rows("4").Delete()
